$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Gioco Cactus / Aggiunta del tutorial / Artegiani): task finished on day 1
$ws.Range("E3").Value = "/"

# Row 4 (Giocatore immune... / Artegiani): remaining effort logged, finished day 3
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "/"

# Row 7 (Quando un giocatore... / Artegiani): remaining effort logged, finished day 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "/"
$ws.Range("F7").Font.ThemeColor = 1
